$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill formulas right from D12:D14 into E:G
$ws.Range("D12:G12").FillRight()
$ws.Range("D13:G13").FillRight()

$ws.Range("E14").Formula = "=IF(OR(ISBLANK(E6), ISBLANK(E7), ISBLANK(E8), ISBLANK(E9), ISBLANK(E10), ISBLANK(E11)),`"`",_xlfn.CONCAT(ROUND(E13*100,0), `"%`"))"
$ws.Range("F14").Formula = "=IF(OR(ISBLANK(F6), ISBLANK(F7), ISBLANK(F8), ISBLANK(F9), ISBLANK(F10), ISBLANK(F11)),`"`",_xlfn.CONCAT(ROUND(F13*100,0), `"%`"))"
$ws.Range("G14").Formula = "=IF(OR(ISBLANK(G6), ISBLANK(G7), ISBLANK(G8), ISBLANK(G9), ISBLANK(G10), ISBLANK(G11)),`"`",_xlfn.CONCAT(ROUND(G13*100,0), `"%`"))"

$ws.Range("B15").Formula = "=A17/(1-A17)"

$ws.Range("G10").Select()
